$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($hf in $sec.Headers) {
        $rng = $hf.Range
        $rng.Find.ClearFormatting()
        $rng.Find.Execute("REPORTE A INVESTIGADORES", $false, $false, $false, $false, $false, $true, 1, $false, "REPORTE A INVESTIGADORES", 2)
    }
}
